$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44964
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 20000
$ws.Range("L2").Value = 21000
$ws.Range("M2").Value = 20500
$ws.Range("P2").Value = 1139

$ws.Range("D3").Value = 45068
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 16000
$ws.Range("L3").Value = 17000
$ws.Range("M3").Value = 16500
$ws.Range("P3").Value = 917

$ws.Range("D4").Value = 44557
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 13000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 13500
$ws.Range("P4").Value = 750

$ws.Range("D5").Value = 44984
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 17000
$ws.Range("L5").Value = 18000
$ws.Range("M5").Value = 17500
$ws.Range("P5").Value = 972

$ws.Range("D6").Value = 44977
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 16500
$ws.Range("L6").Value = 17000
$ws.Range("M6").Value = 16750
$ws.Range("P6").Value = 931

$ws.Range("D7").Value = 44547
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 13000
$ws.Range("L7").Value = 14000
$ws.Range("M7").Value = 13500
$ws.Range("P7").Value = 750

$ws.Range("D8").Value = 45117
$ws.Range("J8").Value = 300
$ws.Range("K8").Value = 17000
$ws.Range("L8").Value = 18000
$ws.Range("M8").Value = 17500
$ws.Range("P8").Value = 972

$ws.Range("D9").Value = 44568
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 15000
$ws.Range("L9").Value = 16000
$ws.Range("M9").Value = 15500
$ws.Range("P9").Value = 861

$ws.Range("D10").Value = 45005
$ws.Range("J10").Value = 200
$ws.Range("K10").Value = 17000
$ws.Range("L10").Value = 18000
$ws.Range("M10").Value = 17500
$ws.Range("P10").Value = 972

$ws.Range("D11").Value = 44998
$ws.Range("J11").Value = 320
$ws.Range("K11").Value = 17000
$ws.Range("L11").Value = 18000
$ws.Range("M11").Value = 17500
$ws.Range("P11").Value = 972

$ws.Range("D12").Value = 44960
$ws.Range("J12").Value = 400
$ws.Range("K12").Value = 19500
$ws.Range("L12").Value = 20000
$ws.Range("M12").Value = 19750
$ws.Range("P12").Value = 1097

$ws.Range("D13").Value = 44957
$ws.Range("J13").Value = 400
$ws.Range("K13").Value = 21000
$ws.Range("L13").Value = 22000
$ws.Range("M13").Value = 21500
$ws.Range("P13").Value = 1194

